$wb = $excel.ActiveWorkbook

# --- "Run Manager" sheet: registrationTest row execute flag Yes -> No ---
$runMgr = $wb.Worksheets.Item("Run Manager")
$runMgr.Range("B3").Value = "No"

# --- "LoginLogoutRegistration" sheet updates ---
$login = $wb.Worksheets.Item("LoginLogoutRegistration")

# browser column D2: chrome -> gdgd
$login.Range("D2").Value = "gdgd"

# execute column C3:C7: Yes -> No
$login.Range("C3").Value = "No"
$login.Range("C4").Value = "No"
$login.Range("C5").Value = "No"
$login.Range("C6").Value = "No"
$login.Range("C7").Value = "No"

# Row 6 registration test data: Moshdada Hamedani -> Boman Irani.
# Leading "'" keeps these as forced-text entries (preserves the existing
# quote-prefixed cell style / shared-string type instead of letting Excel
# re-guess the format, e.g. turning the phone number into a number or the
# email address into an auto-hyperlink-styled cell).
$login.Range("F6").Value = "'Boman"
$login.Range("G6").Value = "'Irani"
$login.Range("H6").Value = "'boman.irani5@gmail.com"
$login.Range("I6").Value = "'6576409987"

# update the active selection to D2 on the active sheet
$login.Activate()
$login.Range("D2").Select()
